$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '67.070.14'
$ws.Cells.Item(2, 5).Value = '  +4.46%  '
$ws.Cells.Item(3, 4).Value = '3.245.90'
$ws.Cells.Item(3, 5).Value = '  +1.97%  '
$ws.Cells.Item(4, 4).Value = "'0.999"
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(4, 5).Value = '  -0.04%  '
$ws.Cells.Item(5, 4).Value = "'577.34"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(6, 4).Value = "'175.86"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  +2.94%  '
$ws.Cells.Item(7, 2).Value = 'XRP'
$ws.Cells.Item(7, 3).Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Cells.Item(7, 4).Value = "'0.606"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = '  +0.53%  '
$ws.Cells.Item(8, 2).Value = 'USDC'
$ws.Cells.Item(8, 3).Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Cells.Item(8, 4).Value = "'0.999"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = '  -0.03%  '
$ws.Cells.Item(9, 4).Value = '3.246.94'
$ws.Cells.Item(9, 5).Value = '  +2.03%  '
$ws.Cells.Item(10, 5).Value = '  +4.76%  '
$ws.Cells.Item(11, 4).Value = "'6.67"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  +1.23%  '
$ws.Cells.Item(12, 4).Value = "'0.406"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  +2.78%  '
$ws.Cells.Item(13, 4).Value = '3.797.16'
$ws.Cells.Item(13, 5).Value = '  +1.82%  '
$ws.Cells.Item(14, 5).Value = '  +2.17%  '
$ws.Cells.Item(15, 4).Value = "'27.75"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = '  +0.99%  '
$ws.Cells.Item(16, 4).Value = '66.961.88'
$ws.Cells.Item(16, 5).Value = '  +3.99%  '
$ws.Cells.Item(17, 5).Value = '  +4.34%  '
$ws.Cells.Item(18, 4).Value = '3.234.37'
$ws.Cells.Item(18, 5).Value = '  +1.71%  '
$ws.Cells.Item(19, 4).Value = "'5.79"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  +2.46%  '
$ws.Cells.Item(20, 4).Value = "'13.21"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  +1.82%  '
$ws.Cells.Item(21, 4).Value = "'366.95"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  +3.93%  '
$ws.Cells.Item(22, 4).Value = "'7.45"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  +3.44%  '
$ws.Cells.Item(23, 5).Value = '  +0.12%  '
$ws.Cells.Item(24, 4).Value = "'69.78"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  +1.66%  '
$ws.Cells.Item(25, 4).Value = "'0.0000119"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  +0.72%  '
$ws.Cells.Item(26, 4).Value = "'0.505"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  +0.49%  '
$ws.Cells.Item(27, 4).Value = '3.363.16'
$ws.Cells.Item(27, 5).Value = '  +1.24%  '
$ws.Cells.Item(28, 4).Value = "'9.83"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  +4.29%  '
$ws.Cells.Item(29, 4).Value = "'0.177"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  +1.61%  '
$ws.Cells.Item(30, 5).Value = '  +0.34%  '
$ws.Cells.Item(31, 5).Value = '  +3.48%  '
$ws.Cells.Item(32, 4).Value = "'5.58"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  -0.43%  '
$ws.Cells.Item(33, 2).Value = 'EthereumClassic'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(33, 4).Value = "'22.42"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  +1.39%  '
$ws.Cells.Item(34, 2).Value = 'USDe'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Cells.Item(34, 4).Value = "'0.998"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  -0.10%  '
$ws.Cells.Item(35, 4).Value = "'1.23"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  +3.27%  '
$ws.Cells.Item(36, 4).Value = "'6.75"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  +2.39%  '
$ws.Cells.Item(37, 4).Value = "'168.35"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  +6.50%  '
$ws.Cells.Item(38, 5).Value = '  +4.33%  '
$ws.Cells.Item(39, 5).Value = '  +3.84%  '
$ws.Cells.Item(40, 4).Value = "'1.85"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  +10.23%  '
$ws.Cells.Item(41, 4).Value = "'26.78"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  +2.07%  '
$ws.Cells.Item(42, 5).Value = '  +4.14%  '
$ws.Cells.Item(43, 4).Value = "'6.34"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  +5.28%  '
$ws.Cells.Item(44, 4).Value = '2.689.21'
$ws.Cells.Item(44, 5).Value = '  +1.48%  '
$ws.Cells.Item(45, 4).Value = "'4.27"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  +3.11%  '
$ws.Cells.Item(46, 4).Value = "'40.40"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  +4.33%  '
$ws.Cells.Item(47, 4).Value = "'0.0672"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  +2.82%  '
$ws.Cells.Item(48, 4).Value = "'24.52"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  +5.37%  '
$ws.Cells.Item(49, 4).Value = "'330.98"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  +3.05%  '
$ws.Cells.Item(50, 4).Value = "'0.0278"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  +3.02%  '
$ws.Cells.Item(51, 4).Value = "'0.103"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  +1.24%  '
